# convenios.xlsx – "Add files via upload" edit
#
# Sheet "Colisiones" (2nd sheet) gains a fourth, right-aligned data column
# ("Máscara que detecta") with per-row values, the header row gets centered,
# and two extra trailing formatted rows appear below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Header row (B2:D2): center-align -------------------------------------
$ws.Range("B2:D2").HorizontalAlignment = -4108   # xlCenter

# --- Row 8 (separator/total row) loses its underline font on column D -----
$ws.Range("D8").Font.Underline = $false

# --- New "Mascara que detecta" column, right-aligned (D3:D18) -------------
$ws.Range("D3:D18").HorizontalAlignment = -4152  # xlRight

$ws.Range("D3").Value = "2,3,5"
$ws.Range("D4").Value = "1,2,5"
$ws.Range("D5").Value = "1,3,5"
$ws.Range("D6").Value = 1.5
$ws.Range("D7").Value = "-"

# --- Leave the cursor where the author's session ended up -----------------
$ws.Range("D8").Select() | Out-Null
